$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H2").Value = 3.2
$ws.Range("I2").Value = 2.63
$ws.Range("M2").Value = 1.04
$ws.Range("O2").Value = 1.25
$ws.Range("Q2").Value = 1.99
$ws.Range("R2").Value = 1.91
$ws.Range("AK2").Value = 23
$ws.Range("AX2").Value = 13
$ws.Range("G3").Value = 2.15
$ws.Range("H3").Value = 2.88
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 3
$ws.Range("M3").Value = 1.1
$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 2.37
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62
$ws.Range("AE3").Value = 19
$ws.Range("AN3").Value = 4
$ws.Range("AO3").Value = 13
$ws.Range("G5").Value = 2.75
$ws.Range("I5").Value = 2.9
$ws.Range("J5").Value = 3.6
$ws.Range("M5").Value = 1.14
$ws.Range("N5").Value = 5.5
$ws.Range("X5").Value = 12
$ws.Range("G9").Value = 1.7
$ws.Range("I9").Value = 4.2
$ws.Range("S9").Value = 1.33
$ws.Range("X9").Value = 9
$ws.Range("AH9").Value = 15
$ws.Range("AJ9").Value = 15
$ws.Range("AL9").Value = 34
$ws.Range("AQ9").Value = 26
$ws.Range("AW9").Value = 6.5
$ws.Range("AX9").Value = 23
$ws.Range("BC9").Value = 501
$ws.Range("S10").Value = 1.5
$ws.Range("M12").Value = 1.05
$ws.Range("N12").Value = 11
$ws.Range("O12").Value = 1.29
$ws.Range("P12").Value = 3.5
$ws.Range("Q12").Value = 1.95
$ws.Range("R12").Value = 1.9
$ws.Range("G14").Value = 2.55
$ws.Range("I14").Value = 2.9
$ws.Range("J14").Value = 3.25
$ws.Range("L14").Value = 3.6
$ws.Range("N14").Value = 8
$ws.Range("W14").Value = 7.5
$ws.Range("X14").Value = 12
$ws.Range("Z14").Value = 26
$ws.Range("AA14").Value = 23
$ws.Range("AH14").Value = 8
$ws.Range("AI14").Value = 13
$ws.Range("AJ14").Value = 11
$ws.Range("AK14").Value = 29
$ws.Range("AL14").Value = 26
$ws.Range("AN14").Value = 4.5
$ws.Range("G17").Value = 3.6
$ws.Range("I17").Value = 2
$ws.Range("X17").Value = 19
$ws.Range("Y17").Value = 13
$ws.Range("AI17").Value = 9.5
$ws.Range("AK17").Value = 17
$ws.Range("AP17").Value = 29
$ws.Range("N20").Value = 8
$ws.Range("G22").Value = 2.88
$ws.Range("H22").Value = 2.8
$ws.Range("I22").Value = 2.75
$ws.Range("L22").Value = 3.4
$ws.Range("M22").Value = 1.1
$ws.Range("N22").Value = 7
$ws.Range("Q22").Value = 2.35
$ws.Range("R22").Value = 1.57
$ws.Range("W22").Value = 7.5
$ws.Range("AA22").Value = 26
$ws.Range("AH22").Value = 7.5
$ws.Range("AI22").Value = 12
$ws.Range("AN22").Value = 4.75
$ws.Range("AO22").Value = 17
$ws.Range("M24").Value = 1.05
$ws.Range("N24").Value = 11
$ws.Range("Q24").Value = 1.98
$ws.Range("R24").Value = 1.88
$ws.Range("G25").Value = 4.33
$ws.Range("H25").Value = 3.75
$ws.Range("I25").Value = 1.75
$ws.Range("J25").Value = 4.5
$ws.Range("L25").Value = 2.3
$ws.Range("Q25").Value = 1.67
$ws.Range("R25").Value = 2.15
$ws.Range("AX25").Value = 9
